# (update) menu pasca penindakan tampilan form create
#
# The only substantive change is the page size of the single section:
#   pgSz w:w="12240" w:h="15840"  (Letter, 8.5in x 11in)
#     -> pgSz w:w="11907" w:h="18711" (Folio/F4, 21cm x 33cm)
#
# Word's PageSetup.PageWidth / PageSetup.PageHeight are expressed in
# points (1 pt = 20 twips), so convert the target twip values to points.

$d = $word.ActiveDocument

$d.PageSetup.PageWidth  = 11907 / 20   # 595.35 pt
$d.PageSetup.PageHeight = 18711 / 20   # 935.55 pt
